$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.323241949081421
$ws.Range("B1").Value = 6.749982833862305
$ws.Range("C1").Value = 6.724996566772461
$ws.Range("D1").Value = 6.755799770355225
$ws.Range("E1").Value = 3.457194089889526
